$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 25 mirrors the date of row 24 (same day, 2019-03-28) and records
# a new practice entry: "647 dp", marked as done.
# Copy A24's formatting (date number format) down into A25 first so the new
# date cell matches the style used by the rest of the date column.
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)

$ws.Range("A25").Value = 43552
$ws.Range("B25").Value = "647 dp"
$ws.Range("F25").Value = "done"

# Leave the selection where the author ended up after entering the data.
$ws.Range("G23").Select()
